$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.210.16"
$ws.Range("E2").Value = "  +0.15%  "
$ws.Range("D3").Value = "1.852.54"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.6985"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.68%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "237.59"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.06%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.0000"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07882"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3014"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.75%  "
$ws.Range("E10").Value = "  +2.76%  "
$ws.Range("E11").Value = "  +0.63%  "
$ws.Range("D12").Value = "1.846.84"
$ws.Range("E12").Value = "  -0.42%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.185"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.02%  "
$ws.Range("E14").Value = "  -2.05%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "89.44"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.28%  "
$ws.Range("D16").Value = "29.210.33"
$ws.Range("E16").Value = "  +0.13%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.800"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.21%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007830"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.34%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.22"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.26%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "235.72"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.73%  "
$ws.Range("E21").Value = "  -0.02%  "
$ws.Range("D22").Value = "2.096.12"
$ws.Range("E22").Value = "  -0.20%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.000"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.08%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.495"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.24%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "162.64"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.67%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.871"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1416"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.04"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.18%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.918"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.75%  "
$ws.Range("E30").Value = "  -0.27%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.473"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.78%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.304"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.26%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.008"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.10%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05146"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.166"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7067"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.74%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9967"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.94%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.677"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.02%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01845"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.12%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.702"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.35%  "
$ws.Range("D41").Value = "1.154.92"
$ws.Range("E41").Value = "  +4.88%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9226"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.19%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.955"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.48%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4236"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.96%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "70.07"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.55%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9998"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.02%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "102.99"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.52%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5293"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.87%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.737"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.06%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.148"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.03%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.951"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.57%  "
